# KIBON-1462 Institutionen Report ergaenzt, kleine Korrekturen
#
# Adds two new columns to the "Data" sheet of the Institutionen report:
#   - "Oeffnungstage" (+ title) right after the "Url" column
#   - "Oeffnungsabweichungen" (+ title) right after the "Oeffnungszeiten" column
#
# Existing columns (and their widths/styles) shift to the right automatically
# via EntireColumn inserts, which also preserve the formatting of the
# surrounding cells (header row style s=6, data row style s=3, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Insert the "Oeffnungstage" column right after "Url" (old column K) ---
$ws.Columns("K").Insert()

# --- Insert the "Oeffnungsabweichungen" column right after "Oeffnungszeiten"
#     (old column K, now shifted to L) ---
$ws.Columns("M").Insert()

# --- Set the new placeholder cells. Order matches how the shared-string
#     table is populated (Title, value, value, Title) so the resulting
#     sharedStrings.xml lines up with the reference workbook. ---
$ws.Range("K4").Value = "{oeffnungstageTitle}"
$ws.Range("K5").Value = "{oeffnungstage}"
$ws.Range("M5").Value = "{oeffnungsAbweichungen}"
$ws.Range("M4").Value = "{oeffnungsAbweichungenTitle}"

# --- Column widths for the two new columns ---
$ws.Columns("K").ColumnWidth = 32.67
$ws.Columns("M").ColumnWidth = 19.33
